$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("@prefix")
$ws.Range("A1").Value = "ome"
$ws = $wb.Worksheets.Item("Image")
$ws.Range("E3").Value = "ome:pixels"
$ws.Range("F3").Value = "ome:acquisitionDate"
$ws.Range("B4").Value = "ome:Image"
$ws.Range("E4").Value = "ome:Pixels"
$ws = $wb.Worksheets.Item("Pixels")
$ws.Range("D3").Value = "ome:pixelType"
$ws.Range("E3").Value = "ome:dimensionOrder"
$ws.Range("F3").Value = "ome:sizeC"
$ws.Range("G3").Value = "ome:sizeT"
$ws.Range("H3").Value = "ome:sizeX"
$ws.Range("I3").Value = "ome:sizeY"
$ws.Range("J3").Value = "ome:sizeZ"
$ws.Range("K3").Value = "ome:channel"
$ws.Range("L3").Value = "ome:binData"
$ws.Range("B4").Value = "ome:Pixels"
$ws.Range("D4").Value = "ome:PixelType"
$ws.Range("E4").Value = "ome:DimensionOrder"
$ws.Range("K4").Value = "ome:Channel"
$ws.Range("L4").Value = "ome:BinData"
$ws = $wb.Worksheets.Item("Channel")
$ws.Range("D3").Value = "ome:color"
$ws.Range("B4").Value = "ome:Channel"
$ws.Range("D4").Value = "ome:Color"
$ws = $wb.Worksheets.Item("Color")
$ws.Range("B4").Value = "ome:Color"
$ws = $wb.Worksheets.Item("Binary_Data")
$ws.Range("C3").Value = "ome:bigEndian"
$ws.Range("D3").Value = "ome:data"
$ws.Range("E3").Value = "ome:length"
$ws.Range("B4").Value = "ome:BinData"
